$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "FirstCene" column (column AG / 33) ---------------------
$ws.Range("AG2").Value = "FirstCene"
for ($r = 3; $r -le 28; $r++) {
    $ws.Cells.Item($r, 33).Value = 0
}
# Row 11 (Square desk) is flagged as part of the first scene
$ws.Range("AG11").Value = 1

# Give the new column (and the one after it) the widths the author left
# them at while working on this area of the sheet.
$ws.Columns.Item(33).ColumnWidth = 9.7109375
$ws.Columns.Item(34).ColumnWidth = 12.140625

# --- Fix up the "Type" column for the desk / mark rows -------------------
$ws.Range("C10").Value = "desk"
$ws.Range("C11").Value = "desk"
$ws.Range("C7").Value = "mark"

# --- Restore the view state the workbook was left in ----------------------
$ws.Range("C12").Select() | Out-Null
